$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'29.501.41"
$ws.Range("E2").Value = "  +1.25%  "

# Row 3
$ws.Range("D3").Value = "'1.923.03"
$ws.Range("E3").Value = "  +2.00%  "

# Row 4
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.44%  "

# Row 5
$ws.Range("D5").Value = "'325.87"
$ws.Range("E5").Value = "  +1.16%  "

# Row 6
$ws.Range("D6").Value = "'1.006"
$ws.Range("E6").Value = "  +0.47%  "

# Row 7
$ws.Range("D7").Value = "'0.4847"
$ws.Range("E7").Value = "  +3.07%  "

# Row 8
$ws.Range("D8").Value = "'0.4104"
$ws.Range("E8").Value = "  +1.94%  "

# Row 9
$ws.Range("D9").Value = "'0.08186"
$ws.Range("E9").Value = "  +2.32%  "

# Row 10
$ws.Range("D10").Value = "'1.026"
$ws.Range("E10").Value = "  +3.41%  "

# Row 11
$ws.Range("D11").Value = "'23.61"
$ws.Range("E11").Value = "  +5.41%  "

# Row 12
$ws.Range("D12").Value = "'1.943.93"
$ws.Range("E12").Value = "  +2.41%  "

# Row 13
$ws.Range("D13").Value = "'6.057"
$ws.Range("E13").Value = "  +3.30%  "

# Row 14
$ws.Range("D14").Value = "'7.255"
$ws.Range("E14").Value = "  +3.39%  "

# Row 15
$ws.Range("D15").Value = "'91.51"
$ws.Range("E15").Value = "  +3.15%  "

# Row 16
$ws.Range("D16").Value = "'0.06783"
$ws.Range("E16").Value = "  +2.51%  "

# Row 17
$ws.Range("E17").Value = "  +0.58%  "

# Row 18
$ws.Range("E18").Value = "  +1.52%  "

# Row 19
$ws.Range("E19").Value = "  +2.58%  "

# Row 20
$ws.Range("E20").Value = "  +0.40%  "

# Row 21
$ws.Range("D21").Value = "'29.535.61"
$ws.Range("E21").Value = "  +1.33%  "

# Row 22
$ws.Range("D22").Value = "'5.637"
$ws.Range("E22").Value = "  +2.85%  "

# Row 23
$ws.Range("D23").Value = "'11.78"
$ws.Range("E23").Value = "  +1.55%  "

# Row 24
$ws.Range("D24").Value = "'2.186"
$ws.Range("E24").Value = "  +0.27%  "

# Row 25
$ws.Range("D25").Value = "'2.120.12"
$ws.Range("E25").Value = "  +0.55%  "

# Row 26
$ws.Range("D26").Value = "'6.739"
$ws.Range("E26").Value = "  +11.86%  "

# Row 27
$ws.Range("D27").Value = "'156.78"
$ws.Range("E27").Value = "  +1.10%  "

# Row 28
$ws.Range("E28").Value = "  +2.70%  "

# Row 29
$ws.Range("D29").Value = "'2.126"
$ws.Range("E29").Value = "  +2.84%  "

# Row 30
$ws.Range("D30").Value = "'120.62"
$ws.Range("E30").Value = "  +2.75%  "

# Row 31
$ws.Range("D31").Value = "'1.033"
$ws.Range("E31").Value = "  -0.09%  "

# Row 32
$ws.Range("D32").Value = "'0.09590"
$ws.Range("E32").Value = "  +1.60%  "

# Row 33
$ws.Range("D33").Value = "'5.536"
$ws.Range("E33").Value = "  +3.66%  "

# Row 34
$ws.Range("D34").Value = "'3.566"
$ws.Range("E34").Value = "  +0.74%  "

# Row 35
$ws.Range("D35").Value = "'1.394"
$ws.Range("E35").Value = "  +1.04%  "

# Row 36
$ws.Range("D36").Value = "'0.02288"
$ws.Range("E36").Value = "  +2.91%  "

# Row 37
$ws.Range("D37").Value = "'0.06150"
$ws.Range("E37").Value = "  +1.46%  "

# Row 38
$ws.Range("D38").Value = "'1.181"
$ws.Range("E38").Value = "  +0.58%  "

# Row 39
$ws.Range("D39").Value = "'0.6000"
$ws.Range("E39").Value = "  +3.45%  "

# Row 40
$ws.Range("D40").Value = "'8.063"
$ws.Range("E40").Value = "  +0.69%  "

# Row 41
$ws.Range("D41").Value = "'10.83"
$ws.Range("E41").Value = "  +8.53%  "

# Row 42
$ws.Range("E42").Value = "  +0.52%  "

# Row 43
$ws.Range("D43").Value = "'0.1867"
$ws.Range("E43").Value = "  +2.44%  "

# Row 44
$ws.Range("D44").Value = "'2.411"
$ws.Range("E44").Value = "  -1.33%  "

# Row 45
$ws.Range("D45").Value = "'1.283"
$ws.Range("E45").Value = "  +1.06%  "

# Row 46
$ws.Range("B46").Value = "Cronos"
$ws.Range("C46").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D46").Value = "'0.07602"
$ws.Range("E46").Value = "  -1.01%  "

# Row 47
$ws.Range("B47").Value = "EnergySwap"
$ws.Range("C47").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D47").Value = "'12.45"
$ws.Range("E47").Value = "  +2.63%  "

# Row 48
$ws.Range("D48").Value = "'0.5598"
$ws.Range("E48").Value = "  +2.49%  "

# Row 49
$ws.Range("D49").Value = "'1.964"
$ws.Range("E49").Value = "  +3.53%  "

# Row 50
$ws.Range("D50").Value = "'117.34"
$ws.Range("E50").Value = "  +3.20%  "

# Row 51
$ws.Range("D51").Value = "'2.439"
$ws.Range("E51").Value = "  +5.09%  "

